# Update the "dSF" column (F) values to repull data / push all data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    4  = -4
    5  = 3
    6  = -3
    7  = 11
    8  = -9
    9  = -4
    12 = 5
    13 = 2
    14 = -1
    15 = -2
    16 = 2
    17 = 6
    18 = -3
    19 = 4
    20 = 1
    21 = -2
    22 = -7
    23 = -2
    24 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
